# display parallel data for tree
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free, explicit cell assignments matching the target diff.
# Rows 8-18: fill in previously-empty Tree BFS/DFS 10000 columns (G, H)
$cells = @{
    "G8"  = 0.1554906
    "H8"  = 0.1367478
    "G9"  = 0.1194432
    "H9"  = 0.1033507
    "A10" = 0.1261244
    "B10" = 0.1748459
    "G10" = 0.1058008
    "H10" = 0.1046368
    "G11" = 0.1483519
    "H11" = 0.1067876
    "G12" = 0.1156435
    "H12" = 0.1086869
    "G13" = 0.1191615
    "H13" = 0.1565017
    "G14" = 0.1511399
    "H14" = 0.1551603
    "G15" = 0.1182255
    "H15" = 0.1031734
    "G16" = 0.1480094
    "H16" = 0.1204223
    "G17" = 0.121699
    "H17" = 0.1217069
    "G18" = 0.1119068
    "H18" = 0.0974884

    # Rows 21-25: fill in previously-empty Graph BFS/DFS 1000 (C, D) and
    # Tree BFS/DFS 50000 (K, L) columns, plus Tree BFS/DFS 1000 (I, J) on 24/25
    "C21" = 0.0100352
    "D21" = 0.014062
    "K21" = 0.5773175
    "L21" = 0.5721709

    "C22" = 0.0177504
    "D22" = 0.0651099
    "K22" = 0.6053008
    "L22" = 0.5751602

    "C23" = 0.0138604
    "D23" = 0.0138828
    "K23" = 0.6084641
    "L23" = 0.6242122

    "I24" = 0.017373
    "J24" = 0.0231225
    "K24" = 0.5986633
    "L24" = 0.5707028

    "I25" = 0.0106078
    "J25" = 0.0140251
    "K25" = 0.547783
    "L25" = 0.5299607

    # New rows 26-27 extend the sheet's data (Tree BFS/DFS 50000 columns)
    "K26" = 0.5881297
    "L26" = 0.6709372

    "K27" = 0.6121012
    "L27" = 0.6250793
}

foreach ($addr in $cells.Keys) {
    $range = $ws.Range($addr)
    $range.ClearFormats()
    $range.Value = $cells[$addr]
}
